$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking "Price" values stay as text (matches source format),
# by temporarily forcing a text number format, then resetting the style so
# no persistent formatting change is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.658.55"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").Value = "3.938.48"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "528.87"
$ws.Range("E5").Value = "  +8.41%  "

$ws.Range("D6").Value = "146.59"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.48%  "

$ws.Range("D10").Value = "0.175"
$ws.Range("E10").Value = "  +5.54%  "

$ws.Range("D11").Value = "0.0000344"
$ws.Range("E11").Value = "  -0.01%  "

$ws.Range("D12").Value = "42.88"
$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "10.51"
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").Value = "4.573.62"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").Value = "3.945.02"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("D16").Value = "14.19"
$ws.Range("E16").Value = "  +0.39%  "

$ws.Range("E17").Value = "  -0.08%  "

$ws.Range("E18").Value = "  +6.75%  "

$ws.Range("D19").Value = "19.90"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("D20").Value = "69.603.96"
$ws.Range("E20").Value = "  +2.00%  "

$ws.Range("D21").Value = "435.51"
$ws.Range("E21").Value = "  +1.14%  "

$ws.Range("E22").Value = "  -4.64%  "

$ws.Range("D23").Value = "14.58"
$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("D24").Value = "88.44"
$ws.Range("E24").Value = "  +1.24%  "

$ws.Range("E25").Value = "  +11.82%  "

$ws.Range("D26").Value = "11.91"
$ws.Range("E26").Value = "  +4.87%  "

$ws.Range("D27").Value = "11.06"
$ws.Range("E27").Value = "  -2.38%  "

$ws.Range("D28").Value = "36.78"
$ws.Range("E28").Value = "  -3.56%  "

$ws.Range("E29").Value = "  -1.33%  "

$ws.Range("D30").Value = "707.79"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").Value = "13.39"
$ws.Range("E31").Value = "  -3.31%  "

$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  -1.41%  "

$ws.Range("D34").Value = "68.19"
$ws.Range("E34").Value = "  +13.13%  "

$ws.Range("E35").Value = "  +8.66%  "

$ws.Range("D36").Value = "6.10"
$ws.Range("E36").Value = "  -2.91%  "

$ws.Range("D37").Value = "0.0₃0872"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "40.54"
$ws.Range("E38").Value = "  -3.12%  "

$ws.Range("E39").Value = "  +0.65%  "

$ws.Range("E40").Value = "  +0.05%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").Value = "0.0486"
$ws.Range("E42").Value = "  +1.55%  "

$ws.Range("E43").Value = "  -3.82%  "

$ws.Range("D44").Value = "3.10"
$ws.Range("E44").Value = "  +6.53%  "

$ws.Range("E45").Value = "  -4.41%  "

$ws.Range("D46").Value = "3.22"
$ws.Range("E46").Value = "  +14.60%  "

$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("E48").Value = "  +1.51%  "

$ws.Range("D49").Value = "0.0₆0368"
$ws.Range("E49").Value = "  +9.21%  "

$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("E51").Value = "  -0.65%  "

# Restore the default style on column D so no stray formatting remains.
$ws.Range("D2:D51").Style = "Normal"
